# Roster ideas: separate "Grupo" (group) block-grade computation out onto
# Sheet2 (previously just a stub "Roster" label), leaving Sheet1's roster
# definition table untouched aside from becoming the non-active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet2: build the "Grupo" block-weighted grade table (B2:J14), a
# scratch weight-from-percentage helper area (N:U), some free-form notes
# (B17:D20), and a merged reminder note (G22:L22).
# ---------------------------------------------------------------------

$ws2.Range("B2").Value = "Roster"

# --- header row (3) ---
$ws2.Range("B3").Value = "Bloque"
$ws2.Range("C3").Value = "w"
$ws2.Range("D3").Value = "wq"
$ws2.Range("E3").Value = "gradeQ"
$ws2.Range("F3").Value = "wt"
$ws2.Range("G3").Value = "gradeT"
$ws2.Range("H3").Value = "wTotal"
$ws2.Range("I3").Value = "finalGrade"
$ws2.Range("J3").Value = "FinalGradeW"

# --- data rows (4-13): Bloque 1..10 ---
$ws2.Range("A4").Value = "Bloque"
$ws2.Range("B4").Value = 1
$ws2.Range("B5").Value = 2
$ws2.Range("B6").Value = 3
$ws2.Range("B7").Value = 4
$ws2.Range("B8").Value = 5
$ws2.Range("B9").Value = 6
$ws2.Range("B10").Value = 7
$ws2.Range("B11").Value = 8
$ws2.Range("B12").Value = 9
$ws2.Range("B13").Value = 10

$ws2.Range("C4").Value = 0
$ws2.Range("C5").Value = 0
$ws2.Range("C6").Value = 0
$ws2.Range("C7").Value = 0
$ws2.Range("C8").Value = 2
$ws2.Range("C9").Value = 2
$ws2.Range("C10").Value = 0
$ws2.Range("C11").Value = 0
$ws2.Range("C12").Value = 0
$ws2.Range("C13").Value = 3

$ws2.Range("D4").Value = 0
$ws2.Range("D5").Value = 0
$ws2.Range("D6").Value = 0
$ws2.Range("D7").Value = 0
$ws2.Range("D8").Value = 3
$ws2.Range("D9").Value = 0
$ws2.Range("D10").Value = 0
$ws2.Range("D11").Value = 0
$ws2.Range("D12").Value = 0
$ws2.Range("D13").Value = 1

$ws2.Range("E4").Value = 0
$ws2.Range("E5").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("E7").Value = 0
$ws2.Range("E8").Value = 100
$ws2.Range("E9").Value = 0
$ws2.Range("E10").Value = 0
$ws2.Range("E11").Value = 0
$ws2.Range("E12").Value = 0
$ws2.Range("E13").Value = 100

$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("F6").Value = 0
$ws2.Range("F7").Value = 0
$ws2.Range("F8").Value = 1
$ws2.Range("F9").Value = 1
$ws2.Range("F10").Value = 0
$ws2.Range("F11").Value = 0
$ws2.Range("F12").Value = 0
$ws2.Range("F13").Value = 1

$ws2.Range("G4").Value = 0
$ws2.Range("G5").Value = 0
$ws2.Range("G6").Value = 0
$ws2.Range("G7").Value = 0
$ws2.Range("G8").Value = 80
$ws2.Range("G9").Value = 80
$ws2.Range("G10").Value = 0
$ws2.Range("G11").Value = 0
$ws2.Range("G12").Value = 0
$ws2.Range("G13").Value = 90

$ws2.Range("H4:H13").Formula = "=D4+F4"
$ws2.Range("I4:I13").Formula = "=IF(H4 > 0,((D4*E4)+(F4*G4))/H4,0)"
$ws2.Range("J4:J13").Formula = "=IF(C4>0,C4*I4,0)"

# --- totals row (14) ---
$ws2.Range("A14").Value = "Grupo"
$ws2.Range("B14").Value = "FG"
$ws2.Range("C14").Formula = "=SUM(C4:C13)"
$ws2.Range("D14").Value = $null
$ws2.Range("E14").Value = $null
$ws2.Range("F14").Value = $null
$ws2.Range("G14").Value = $null
$ws2.Range("H14").Value = $null
$ws2.Range("I14").Formula = "=SUM(J4:J13)/C14"
